$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "D2" = "29.078.04"
    "E2" = "  +0.15%  "
    "D3" = "1.835.22"
    "E3" = "  +0.23%  "
    "D4" = "0.9996"
    "E4" = "  +0.14%  "
    "D5" = "242.84"
    "E5" = "  -0.71%  "
    "D6" = "0.6140"
    "E6" = "  -2.87%  "
    "E7" = "  +0.20%  "
    "D8" = "0.07465"
    "E8" = "  -0.67%  "
    "D9" = "0.2918"
    "E9" = "  -0.67%  "
    "D10" = "23.13"
    "E10" = "  +1.05%  "
    "D11" = "0.07686"
    "E11" = "  -0.20%  "
    "D12" = "1.842.02"
    "E12" = "  +0.78%  "
    "D13" = "5.005"
    "E13" = "  +0.20%  "
    "D14" = "0.6718"
    "E14" = "  +0.22%  "
    "D15" = "82.66"
    "E15" = "  -0.41%  "
    "D16" = "0.000009146"
    "E16" = "  -4.61%  "
    "D17" = "5.925"
    "E17" = "  -2.53%  "
    "D18" = "29.074.44"
    "D19" = "2.082.12"
    "E19" = "  +0.28%  "
    "D20" = "231.70"
    "E20" = "  +2.32%  "
    "E21" = "  +0.58%  "
    "E22" = "  +0.29%  "
    "D23" = "7.196"
    "E23" = "  +0.53%  "
    "D24" = "1.001"
    "E24" = "  +0.15%  "
    "D25" = "159.48"
    "E25" = "  -0.39%  "
    "D26" = "0.1387"
    "E26" = "  -1.42%  "
    "D27" = "8.496"
    "E27" = "  -0.49%  "
    "E28" = "  -0.69%  "
    "D29" = "1.493"
    "E29" = "  -0.37%  "
    "E30" = "  +0.70%  "
    "D31" = "4.137"
    "E31" = "  +1.77%  "
    "D32" = "0.05549"
    "E32" = "  +3.10%  "
    "D33" = "1.208"
    "E33" = "  +0.76%  "
    "D34" = "0.7443"
    "E34" = "  +0.02%  "
    "D35" = "1.836"
    "E35" = "  -1.20%  "
    "D36" = "1.140"
    "E36" = "  +0.08%  "
    "D37" = "2.660"
    "E37" = "  +0.18%  "
    "D38" = "2.770"
    "E38" = "  +0.15%  "
    "E39" = "  -0.65%  "
    "D40" = "1.209.00"
    "E40" = "  -2.78%  "
    "D41" = "6.480"
    "E41" = "  -2.50%  "
    "D42" = "0.8931"
    "E42" = "  -1.16%  "
    "E43" = "  +0.06%  "
    "D44" = "101.99"
    "E44" = "  +0.21%  "
    "D45" = "1.984.44"
    "D46" = "65.66"
    "E46" = "  +1.29%  "
    "B47" = "Mantle"
    "C47" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "D47" = "0.5089"
    "E47" = "  -0.32%  "
    "B48" = "TheSandbox"
    "C48" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
    "D48" = "0.4066"
    "E48" = "  -0.15%  "
    "B49" = "EnergySwap"
    "C49" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D49" = "9.108"
    "E49" = "  +1.16%  "
    "B50" = "Cronos"
    "C50" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D50" = "0.05830"
    "E50" = "  +1.05%  "
    "B51" = "XinFinNetwork"
    "C51" = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
    "D51" = "0.07173"
    "E51" = "  +8.37%  "
}

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
}
